# Update the "Corr/total marks" section of the marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: Marking - Right value changes from 3 to 5
$ws.Range("B11").Value = 5

# B12: Total - Right value changes from 54 to 90
$ws.Range("B12").Value = 90

# E12: Total - Max summary text changes from "50/84" to "90/140"
$ws.Range("E12").Value = "90/140"
